# edit.ps1
# Applies the "update scripts wuth new tpm" commit to Fgf1-Fgfr2.xlsx
#
# Summary of the change:
#  1. Cluster label rename (affects the shared-string table, which in turn
#     changes the displayed "Sending cluster" (col A) / "Target cluster" (col D)
#     text for every row that referenced the renamed labels):
#        "Inflammatory-Mac" -> "MuSCs"
#        "MuSCs"            -> "Neutrophils"
#  2. Recomputed TPM-derived NATMI statistics (columns E-T) for every data row
#     (rows 2-16), reflecting the new underlying TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the two cluster labels everywhere they are used -------------
# Rows 8-10 were sent from "Inflammatory-Mac" -> now "MuSCs"
# Rows 11-13 were sent from "MuSCs" -> now "Neutrophils"
# (Target-cluster column D only ever referenced "MuSCs" in rows 4/7/10/13/16,
#  those are rewritten below together with their numeric updates.)
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(11,1).Value = "Neutrophils"
$ws.Cells.Item(12,1).Value = "Neutrophils"
$ws.Cells.Item(13,1).Value = "Neutrophils"

# --- 2) Updated TPM-derived values (and the D-column cluster renames) ------
    # Row 2
    $ws.Cells.Item(2,7).Value = 0.6874376666666667
    $ws.Cells.Item(2,8).Value = 2.062313
    $ws.Cells.Item(2,9).Value = 0.2330845252991127
    $ws.Cells.Item(2,10).Value = 0.2330845252991127
    $ws.Cells.Item(2,13).Value = 0.106124
    $ws.Cells.Item(2,14).Value = 0.318372
    $ws.Cells.Item(2,15).Value = 0.08094716512538251
    $ws.Cells.Item(2,16).Value = 0.08094716512538253
    $ws.Cells.Item(2,17).Value = 0.07295363493733334
    $ws.Cells.Item(2,18).Value = 0.656582714436
    $ws.Cells.Item(2,19).Value = 0.01886753155755868
    $ws.Cells.Item(2,20).Value = 0.01886753155755868
    # Row 3
    $ws.Cells.Item(3,7).Value = 0.6874376666666667
    $ws.Cells.Item(3,8).Value = 2.062313
    $ws.Cells.Item(3,9).Value = 0.2330845252991127
    $ws.Cells.Item(3,10).Value = 0.2330845252991127
    $ws.Cells.Item(3,15).Value = 0.8331551016962769
    $ws.Cells.Item(3,16).Value = 0.833155101696277
    $ws.Cells.Item(3,17).Value = 0.7508810597774445
    $ws.Cells.Item(3,18).Value = 6.757929537997001
    $ws.Cells.Item(3,19).Value = 0.1941955613794107
    $ws.Cells.Item(3,20).Value = 0.1941955613794107
    # Row 4
    $ws.Cells.Item(4,4).Value = "MuSCs"
    $ws.Cells.Item(4,7).Value = 0.6874376666666667
    $ws.Cells.Item(4,8).Value = 2.062313
    $ws.Cells.Item(4,9).Value = 0.2330845252991127
    $ws.Cells.Item(4,10).Value = 0.2330845252991127
    $ws.Cells.Item(4,13).Value = 0.1126143333333333
    $ws.Cells.Item(4,14).Value = 0.337843
    $ws.Cells.Item(4,15).Value = 0.08589773317834044
    $ws.Cells.Item(4,16).Value = 0.08589773317834046
    $ws.Cells.Item(4,17).Value = 0.0774153345398889
    $ws.Cells.Item(4,18).Value = 0.696738010859
    $ws.Cells.Item(4,19).Value = 0.02002143236214333
    $ws.Cells.Item(4,20).Value = 0.02002143236214333
    # Row 5
    $ws.Cells.Item(5,9).Value = 0.1587189032810992
    $ws.Cells.Item(5,10).Value = 0.1587189032810992
    $ws.Cells.Item(5,13).Value = 0.106124
    $ws.Cells.Item(5,14).Value = 0.318372
    $ws.Cells.Item(5,15).Value = 0.08094716512538251
    $ws.Cells.Item(5,16).Value = 0.08094716512538253
    $ws.Cells.Item(5,17).Value = 0.04967777638933334
    $ws.Cells.Item(5,18).Value = 0.447099987504
    $ws.Cells.Item(5,19).Value = 0.01284784527241475
    $ws.Cells.Item(5,20).Value = 0.01284784527241475
    # Row 6
    $ws.Cells.Item(6,9).Value = 0.1587189032810992
    $ws.Cells.Item(6,10).Value = 0.1587189032810992
    $ws.Cells.Item(6,15).Value = 0.8331551016962769
    $ws.Cells.Item(6,16).Value = 0.833155101696277
    $ws.Cells.Item(6,19).Value = 0.1322374640042858
    $ws.Cells.Item(6,20).Value = 0.1322374640042858
    # Row 7
    $ws.Cells.Item(7,4).Value = "MuSCs"
    $ws.Cells.Item(7,9).Value = 0.1587189032810992
    $ws.Cells.Item(7,10).Value = 0.1587189032810992
    $ws.Cells.Item(7,13).Value = 0.1126143333333333
    $ws.Cells.Item(7,14).Value = 0.337843
    $ws.Cells.Item(7,15).Value = 0.08589773317834044
    $ws.Cells.Item(7,16).Value = 0.08589773317834046
    $ws.Cells.Item(7,17).Value = 0.0527159706528889
    $ws.Cells.Item(7,18).Value = 0.474443735876
    $ws.Cells.Item(7,19).Value = 0.01363359400439869
    $ws.Cells.Item(7,20).Value = 0.01363359400439869
    # Row 8
    $ws.Cells.Item(8,5).Value = 3
    $ws.Cells.Item(8,6).Value = 1
    $ws.Cells.Item(8,7).Value = 1.758325333333333
    $ws.Cells.Item(8,8).Value = 5.274976
    $ws.Cells.Item(8,9).Value = 0.5961826730104559
    $ws.Cells.Item(8,10).Value = 0.5961826730104558
    $ws.Cells.Item(8,13).Value = 0.106124
    $ws.Cells.Item(8,14).Value = 0.318372
    $ws.Cells.Item(8,15).Value = 0.08094716512538251
    $ws.Cells.Item(8,16).Value = 0.08094716512538253
    $ws.Cells.Item(8,17).Value = 0.1866005176746666
    $ws.Cells.Item(8,18).Value = 1.679404659072
    $ws.Cells.Item(8,19).Value = 0.0482592972770693
    $ws.Cells.Item(8,20).Value = 0.0482592972770693
    # Row 9
    $ws.Cells.Item(9,5).Value = 3
    $ws.Cells.Item(9,6).Value = 1
    $ws.Cells.Item(9,7).Value = 1.758325333333333
    $ws.Cells.Item(9,8).Value = 5.274976
    $ws.Cells.Item(9,9).Value = 0.5961826730104559
    $ws.Cells.Item(9,10).Value = 0.5961826730104558
    $ws.Cells.Item(9,15).Value = 0.8331551016962769
    $ws.Cells.Item(9,16).Value = 0.833155101696277
    $ws.Cells.Item(9,17).Value = 1.920600592238222
    $ws.Cells.Item(9,18).Value = 17.285405330144
    $ws.Cells.Item(9,19).Value = 0.4967126355615846
    $ws.Cells.Item(9,20).Value = 0.4967126355615846
    # Row 10
    $ws.Cells.Item(10,4).Value = "MuSCs"
    $ws.Cells.Item(10,5).Value = 3
    $ws.Cells.Item(10,6).Value = 1
    $ws.Cells.Item(10,7).Value = 1.758325333333333
    $ws.Cells.Item(10,8).Value = 5.274976
    $ws.Cells.Item(10,9).Value = 0.5961826730104559
    $ws.Cells.Item(10,10).Value = 0.5961826730104558
    $ws.Cells.Item(10,13).Value = 0.1126143333333333
    $ws.Cells.Item(10,14).Value = 0.337843
    $ws.Cells.Item(10,15).Value = 0.08589773317834044
    $ws.Cells.Item(10,16).Value = 0.08589773317834046
    $ws.Cells.Item(10,17).Value = 0.1980126351964444
    $ws.Cells.Item(10,18).Value = 1.782113716768
    $ws.Cells.Item(10,19).Value = 0.05121074017180193
    $ws.Cells.Item(10,20).Value = 0.05121074017180193
    # Row 11
    $ws.Cells.Item(11,5).Value = 1
    $ws.Cells.Item(11,6).Value = 0.3333333333333333
    $ws.Cells.Item(11,7).Value = 0.01541033333333333
    $ws.Cells.Item(11,8).Value = 0.046231
    $ws.Cells.Item(11,9).Value = 0.005225070437466708
    $ws.Cells.Item(11,10).Value = 0.005225070437466708
    $ws.Cells.Item(11,13).Value = 0.106124
    $ws.Cells.Item(11,14).Value = 0.318372
    $ws.Cells.Item(11,15).Value = 0.08094716512538251
    $ws.Cells.Item(11,16).Value = 0.08094716512538253
    $ws.Cells.Item(11,17).Value = 0.001635406214666667
    $ws.Cells.Item(11,18).Value = 0.014718655932
    $ws.Cells.Item(11,19).Value = 0.0004229546394933722
    $ws.Cells.Item(11,20).Value = 0.0004229546394933723
    # Row 12
    $ws.Cells.Item(12,5).Value = 1
    $ws.Cells.Item(12,6).Value = 0.3333333333333333
    $ws.Cells.Item(12,7).Value = 0.01541033333333333
    $ws.Cells.Item(12,8).Value = 0.046231
    $ws.Cells.Item(12,9).Value = 0.005225070437466708
    $ws.Cells.Item(12,10).Value = 0.005225070437466708
    $ws.Cells.Item(12,15).Value = 0.8331551016962769
    $ws.Cells.Item(12,16).Value = 0.833155101696277
    $ws.Cells.Item(12,17).Value = 0.01683254785988889
    $ws.Cells.Item(12,18).Value = 0.151492930739
    $ws.Cells.Item(12,19).Value = 0.004353294091697785
    $ws.Cells.Item(12,20).Value = 0.004353294091697785
    # Row 13
    $ws.Cells.Item(13,4).Value = "MuSCs"
    $ws.Cells.Item(13,5).Value = 1
    $ws.Cells.Item(13,6).Value = 0.3333333333333333
    $ws.Cells.Item(13,7).Value = 0.01541033333333333
    $ws.Cells.Item(13,8).Value = 0.046231
    $ws.Cells.Item(13,9).Value = 0.005225070437466708
    $ws.Cells.Item(13,10).Value = 0.005225070437466708
    $ws.Cells.Item(13,13).Value = 0.1126143333333333
    $ws.Cells.Item(13,14).Value = 0.337843
    $ws.Cells.Item(13,15).Value = 0.08589773317834044
    $ws.Cells.Item(13,16).Value = 0.08589773317834046
    $ws.Cells.Item(13,17).Value = 0.001735424414777778
    $ws.Cells.Item(13,18).Value = 0.015618819733
    $ws.Cells.Item(13,19).Value = 0.0004488217062755498
    $ws.Cells.Item(13,20).Value = 0.0004488217062755499
    # Row 14
    $ws.Cells.Item(14,7).Value = 0.02002233333333333
    $ws.Cells.Item(14,8).Value = 0.060067
    $ws.Cells.Item(14,9).Value = 0.006788827971865474
    $ws.Cells.Item(14,10).Value = 0.006788827971865474
    $ws.Cells.Item(14,13).Value = 0.106124
    $ws.Cells.Item(14,14).Value = 0.318372
    $ws.Cells.Item(14,15).Value = 0.08094716512538251
    $ws.Cells.Item(14,16).Value = 0.08094716512538253
    $ws.Cells.Item(14,17).Value = 0.002124850102666667
    $ws.Cells.Item(14,18).Value = 0.019123650924
    $ws.Cells.Item(14,19).Value = 0.0005495363788464102
    $ws.Cells.Item(14,20).Value = 0.0005495363788464103
    # Row 15
    $ws.Cells.Item(15,7).Value = 0.02002233333333333
    $ws.Cells.Item(15,8).Value = 0.060067
    $ws.Cells.Item(15,9).Value = 0.006788827971865474
    $ws.Cells.Item(15,10).Value = 0.006788827971865474
    $ws.Cells.Item(15,15).Value = 0.8331551016962769
    $ws.Cells.Item(15,16).Value = 0.833155101696277
    $ws.Cells.Item(15,17).Value = 0.02187018780255555
    $ws.Cells.Item(15,18).Value = 0.196831690223
    $ws.Cells.Item(15,19).Value = 0.005656146659298109
    $ws.Cells.Item(15,20).Value = 0.005656146659298109
    # Row 16
    $ws.Cells.Item(16,4).Value = "MuSCs"
    $ws.Cells.Item(16,7).Value = 0.02002233333333333
    $ws.Cells.Item(16,8).Value = 0.060067
    $ws.Cells.Item(16,9).Value = 0.006788827971865474
    $ws.Cells.Item(16,10).Value = 0.006788827971865474
    $ws.Cells.Item(16,13).Value = 0.1126143333333333
    $ws.Cells.Item(16,14).Value = 0.337843
    $ws.Cells.Item(16,15).Value = 0.08589773317834044
    $ws.Cells.Item(16,16).Value = 0.08589773317834046
    $ws.Cells.Item(16,17).Value = 0.002254801720111111
    $ws.Cells.Item(16,18).Value = 0.020293215481
    $ws.Cells.Item(16,19).Value = 0.0005831449337209546
    $ws.Cells.Item(16,20).Value = 0.0005831449337209547
